# "fixed 2nd diagram again"
# Wrap the "Arc 33" shape and the "Group 121" group (the 2nd diagram's
# sibling shapes sitting directly on the slide) into a single new group,
# matching PowerPoint's native Group() behavior (new group becomes
# "Group 1", id reset to 2, inserted at the position of the first
# selected shape).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$arcShape = $null
$groupShape = $null

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Arc 33") {
        $arcShape = $shp
    }
    if ($shp.Name -eq "Group 121") {
        $groupShape = $shp
    }
}

$range = $s.Shapes.Range(@($arcShape.Name, $groupShape.Name))
$range.Group() | Out-Null
